$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.166.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.644.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.81%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -1.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.140"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.66%  "
$ws.Range("E10").Value = "  -1.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.25"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.352"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("E13").Value = "  -0.31%  "
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.124.85"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "68.050.22"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.649.51"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.38"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "363.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.05%  "
$ws.Range("E20").Value = "  -1.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.74%  "
$ws.Range("E22").Value = "  -2.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "75.13"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.45%  "
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.89%  "
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("E29").Value = "  -1.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "554.32"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.04"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.41"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.53%  "
$ws.Range("E33").Value = "  -1.04%  "
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("E35").Value = "  -2.34%  "
$ws.Range("E36").Value = "  -0.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "161.31"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.54"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.88%  "
$ws.Range("E39").Value = "  +0.65%  "
$ws.Range("E40").Value = "  -3.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.33"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0₆0334"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.30%  "
$ws.Range("E43").Value = "  +0.24%  "
$ws.Range("E44").Value = "  -2.83%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "158.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.42%  "
$ws.Range("E47").Value = "  -0.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("E49").Value = "  -2.23%  "
$ws.Range("E50").Value = "  -0.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.615"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.96%  "
